# "add results from latest run"
# Updates the Nowcasts table: refreshed numeric results for the existing
# rows (2025-03-30 .. 2025-08-15) and appends a new row for 2025-08-30,
# plus a minor column-F width tweak.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (2025-03-30): revised values -----------------------------------
$ws.Cells.Item(2,2).Value  = [double]"0.29327604952437519"
$ws.Cells.Item(2,3).Value  = [double]"0"
$ws.Cells.Item(2,4).Value  = [double]"0"
$ws.Cells.Item(2,5).Value  = [double]"0"
$ws.Cells.Item(2,6).Value  = [double]"0"
$ws.Cells.Item(2,7).Value  = [double]"0"
$ws.Cells.Item(2,8).Value  = [double]"0"
$ws.Cells.Item(2,9).Value  = [double]"0"
$ws.Cells.Item(2,10).Value = [double]"0"
$ws.Cells.Item(2,11).Value = [double]"0"

# --- Row 3 (2025-04-15): revised values -----------------------------------
$ws.Cells.Item(3,2).Value  = [double]"0.29425048920040919"
$ws.Cells.Item(3,3).Value  = [double]"0"
$ws.Cells.Item(3,4).Value  = [double]"0.0013189885815270894"
$ws.Cells.Item(3,5).Value  = [double]"-4.999241294972401e-06"
$ws.Cells.Item(3,6).Value  = [double]"-3.1859137857299906e-06"
$ws.Cells.Item(3,7).Value  = [double]"-0.00013462783614299903"
$ws.Cells.Item(3,8).Value  = [double]"7.3893907311130055e-07"
$ws.Cells.Item(3,9).Value  = [double]"-0.00020239736080266213"
$ws.Cells.Item(3,10).Value = [double]"0"
$ws.Cells.Item(3,11).Value = [double]"-7.7492539829471241e-08"

# --- Row 4 (2025-04-30): revised values -----------------------------------
$ws.Cells.Item(4,2).Value  = [double]"0.29145801663328447"
$ws.Cells.Item(4,3).Value  = [double]"-0.0015539283330628685"
$ws.Cells.Item(4,4).Value  = [double]"0"
$ws.Cells.Item(4,5).Value  = [double]"-1.4930367743056516e-07"
$ws.Cells.Item(4,6).Value  = [double]"-2.0086193578360719e-05"
$ws.Cells.Item(4,7).Value  = [double]"0"
$ws.Cells.Item(4,8).Value  = [double]"-4.2873833285103199e-05"
$ws.Cells.Item(4,9).Value  = [double]"-0.0011575960473246382"
$ws.Cells.Item(4,10).Value = [double]"2.5560372760684001e-08"
$ws.Cells.Item(4,11).Value = [double]"-1.7864416569091102e-05"

# --- Row 5 (2025-05-15): revised values -----------------------------------
$ws.Cells.Item(5,2).Value  = [double]"0.29203767419653448"
$ws.Cells.Item(5,3).Value  = [double]"0.0032496962659846035"
$ws.Cells.Item(5,4).Value  = [double]"-0.00062363513829834618"
$ws.Cells.Item(5,5).Value  = [double]"6.2383460227012725e-05"
$ws.Cells.Item(5,6).Value  = [double]"0.00016816307388065456"
$ws.Cells.Item(5,7).Value  = [double]"-0.0018508191266549041"
$ws.Cells.Item(5,8).Value  = [double]"5.6719468706137904e-06"
$ws.Cells.Item(5,9).Value  = [double]"-0.00043254868239875891"
$ws.Cells.Item(5,10).Value = [double]"0"
$ws.Cells.Item(5,11).Value = [double]"7.4576363912060017e-07"

# --- Row 6 (2025-05-30): revised values -----------------------------------
$ws.Cells.Item(6,2).Value  = [double]"0.58250535049530683"
$ws.Cells.Item(6,3).Value  = [double]"0.2953383645413527"
$ws.Cells.Item(6,4).Value  = [double]"0"
$ws.Cells.Item(6,5).Value  = [double]"-3.1956150137918815e-07"
$ws.Cells.Item(6,6).Value  = [double]"-8.7279565938449616e-05"
$ws.Cells.Item(6,7).Value  = [double]"0"
$ws.Cells.Item(6,8).Value  = [double]"-5.3246352312610439e-05"
$ws.Cells.Item(6,9).Value  = [double]"-0.0047294465563309223"
$ws.Cells.Item(6,10).Value = [double]"0"
$ws.Cells.Item(6,11).Value = [double]"-3.9620649699978472e-07"

# --- Row 7 (2025-06-15): revised values -----------------------------------
$ws.Cells.Item(7,2).Value  = [double]"0.54039031374823476"
$ws.Cells.Item(7,3).Value  = [double]"0"
$ws.Cells.Item(7,4).Value  = [double]"-0.043212493197749016"
$ws.Cells.Item(7,5).Value  = [double]"-9.6715138137333262e-05"
$ws.Cells.Item(7,6).Value  = [double]"-0.00077516748501435572"
$ws.Cells.Item(7,7).Value  = [double]"0.00033905404093473584"
$ws.Cells.Item(7,8).Value  = [double]"0"
$ws.Cells.Item(7,9).Value  = [double]"0.0016385504979685793"
$ws.Cells.Item(7,10).Value = [double]"0"
$ws.Cells.Item(7,11).Value = [double]"-8.265465074708267e-06"

# --- Row 8 (2025-06-30): revised values -----------------------------------
$ws.Cells.Item(8,2).Value  = [double]"0.23692802112716793"
$ws.Cells.Item(8,3).Value  = [double]"-0.29069282874203783"
$ws.Cells.Item(8,4).Value  = [double]"0"
$ws.Cells.Item(8,5).Value  = [double]"2.0139816837635527e-05"
$ws.Cells.Item(8,6).Value  = [double]"-0.007431737238923701"
$ws.Cells.Item(8,7).Value  = [double]"0"
$ws.Cells.Item(8,8).Value  = [double]"-2.7217076120107353e-05"
$ws.Cells.Item(8,9).Value  = [double]"-0.0054883165996474723"
$ws.Cells.Item(8,10).Value = [double]"0"
$ws.Cells.Item(8,11).Value = [double]"0.00015766721882470858"

# --- Row 9 (2025-07-15): revised values -----------------------------------
$ws.Cells.Item(9,2).Value  = [double]"-0.015116105433552685"
$ws.Cells.Item(9,3).Value  = [double]"0"
$ws.Cells.Item(9,4).Value  = [double]"-0.031614122423793321"
$ws.Cells.Item(9,5).Value  = [double]"-0.029457147468180633"
$ws.Cells.Item(9,6).Value  = [double]"-0.18359885899239659"
$ws.Cells.Item(9,7).Value  = [double]"-0.0027840390600413676"
$ws.Cells.Item(9,8).Value  = [double]"-0.0046598565985617206"
$ws.Cells.Item(9,9).Value  = [double]"-0.00023216564908855146"
$ws.Cells.Item(9,10).Value = [double]"0"
$ws.Cells.Item(9,11).Value = [double]"0.00030206363134155367"

# --- Row 10 (2025-07-30): revised values ----------------------------------
$ws.Cells.Item(10,2).Value  = [double]"0.28640484474548183"
$ws.Cells.Item(10,3).Value  = [double]"0.33611848773701059"
$ws.Cells.Item(10,4).Value  = [double]"0"
$ws.Cells.Item(10,5).Value  = [double]"-0.00062088709854028539"
$ws.Cells.Item(10,6).Value  = [double]"0.0063813633989069276"
$ws.Cells.Item(10,7).Value  = [double]"0"
$ws.Cells.Item(10,8).Value  = [double]"-0.0003456723929400724"
$ws.Cells.Item(10,9).Value  = [double]"-0.00086665215979439959"
$ws.Cells.Item(10,10).Value = [double]"-0.048497421483609159"
$ws.Cells.Item(10,11).Value = [double]"0.0093517321780009488"

# --- Row 11 (2025-08-15): revised values ----------------------------------
$ws.Cells.Item(11,2).Value  = [double]"0.42001279568163252"
$ws.Cells.Item(11,3).Value  = [double]"0"
$ws.Cells.Item(11,4).Value  = [double]"-0.012226240649239592"
$ws.Cells.Item(11,5).Value  = [double]"0.022779814260864957"
$ws.Cells.Item(11,6).Value  = [double]"0.16181858389951553"
$ws.Cells.Item(11,7).Value  = [double]"0.0015996221714689408"
$ws.Cells.Item(11,8).Value  = [double]"8.677072995056951e-06"
$ws.Cells.Item(11,9).Value  = [double]"-0.010560155435089253"
$ws.Cells.Item(11,10).Value = [double]"0"
$ws.Cells.Item(11,11).Value = [double]"-0.029812350384364972"

# --- Row 12 (new, 2025-08-30) ----------------------------------------------
# Force text so the date-like label isn't auto-converted to a date serial,
# then drop back to the workbook's default "Normal" style so the cell
# matches the plain (unstyled) look of the other row labels.
$ws.Cells.Item(12,1).NumberFormat = "@"
$ws.Cells.Item(12,1).Value = "2025-08-30"
$ws.Cells.Item(12,1).Style = "Normal"

$ws.Cells.Item(12,2).Value  = [double]"0.32696495078983284"
$ws.Cells.Item(12,3).Value  = [double]"-0.042561976288722506"
$ws.Cells.Item(12,4).Value  = [double]"0"
$ws.Cells.Item(12,5).Value  = [double]"0.0011715483887053417"
$ws.Cells.Item(12,6).Value  = [double]"9.676323744577445e-05"
$ws.Cells.Item(12,7).Value  = [double]"0"
$ws.Cells.Item(12,8).Value  = [double]"3.4774571382049886e-05"
$ws.Cells.Item(12,9).Value  = [double]"-0.0076912114785374856"
$ws.Cells.Item(12,10).Value = [double]"0"
$ws.Cells.Item(12,11).Value = [double]"-0.044097743322072835"

# --- Column F width: 16.24609375 -> 15.77734375 ----------------------------
$ws.Columns.Item(6).ColumnWidth = 15
